$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new row 2 for "Wassertunnel, Wasserstollen, Druckstollen" / 2023 / caveBuilding ---
$ws.Rows.Item(2).Insert()

$ws.Range("A2").Value = "Wassertunnel, Wasserstollen, Druckstollen"
$ws.Range("B2").Value = 2023
$ws.Range("A2").Font.Bold = $false
$ws.Range("B2").Font.Bold = $false

$ws.Range("C2").Value = "http://inspire.ec.europa.eu/codelist/BuildingNatureValue/caveBuilding"
$ws.Hyperlinks.Add($ws.Range("C2"), "http://inspire.ec.europa.eu/codelist/BuildingNatureValue/caveBuilding")

# --- Insert a new row 9 (after "Sperrwerk", before "Schöpfwerk") for "Verschlussbauwerk" / 2085 / dam ---
$ws.Rows.Item(9).Insert()

$ws.Range("A9").Value = "Verschlussbauwerk"
$ws.Range("B9").Value = 2085
$ws.Range("C9").Value = "http://inspire.ec.europa.eu/codelist/BuildingNatureValue/dam"
$ws.Range("C9").Style = "Link"

# --- Restore active selection to C9, matching the edited workbook's cursor position ---
$ws.Range("C9").Select()
